$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.500.66"
$ws.Range("E2").Value = "  +3.93%  "
$ws.Range("D3").Value = "3.274.84"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'213.33"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").Value = "'629.62"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +21.99%  "
$ws.Range("D8").Value = "'0.697"
$ws.Range("E8").Value = "  +17.63%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "3.270.61"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "'0.578"
$ws.Range("E11").Value = "  -3.28%  "
$ws.Range("E12").Value = "  +10.46%  "
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").Value = "'34.28"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "3.876.11"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "'5.35"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "88.254.87"
$ws.Range("E17").Value = "  +4.03%  "
$ws.Range("D18").Value = "3.282.27"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").Value = "'3.13"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'14.15"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "'436.04"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "'8.91"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'7.36"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "'12.37"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("D26").Value = "'5.19"
$ws.Range("E26").Value = "  -5.24%  "
$ws.Range("D27").Value = "3.434.60"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "'77.07"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'0.185"
$ws.Range("E31").Value = "  +13.52%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'8.90"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'571.45"
$ws.Range("E34").Value = "  -5.65%  "
$ws.Range("D35").Value = "'1.40"
$ws.Range("E35").Value = "  -9.87%  "
$ws.Range("E36").Value = "  -3.82%  "
$ws.Range("D37").Value = "'7.11"
$ws.Range("E37").Value = "  +11.02%  "
$ws.Range("E38").Value = "  -8.29%  "
$ws.Range("D39").Value = "'22.66"
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'21.78"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("D43").Value = "'2.03"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "'2.98"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'153.17"
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("D47").Value = "'180.31"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").Value = "'44.79"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "'1.29"
$ws.Range("E49").Value = "  -4.74%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").Value = "'0.0680"
$ws.Range("E50").Value = "  +21.76%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.24"
$ws.Range("E51").Value = "  -0.20%  "
